$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 58 (shifts existing rows 58:86 down to 59:87)
$ws.Rows("58:58").Insert()

# Populate the new row 58 with the weekly Mango price update
$ws.Range("A58").Value = 1
$ws.Range("B58").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C58").Value = "Arica y Parinacota"
$ws.Range("D58").Value = 44447
$ws.Range("E58").Value = 15
$ws.Range("F58").Value = "Fruta"
$ws.Range("G58").Value = 100108
$ws.Range("H58").Value = "Tropicales y subtropicales"
$ws.Range("I58").Value = 100108002
$ws.Range("J58").Value = "Mango"
$ws.Range("K58").Value = "Sin especificar"
$ws.Range("L58").Value = "Extra"
$ws.Range("M58").Value = 300
$ws.Range("N58").Value = 6800
$ws.Range("O58").Value = 7000
$ws.Range("P58").Value = 6900
$ws.Range("Q58").Value = "$/bandeja 4 kilos"
$ws.Range("R58").Value = "Perú"
$ws.Range("S58").Value = 1725
$ws.Range("T58").Value = 4
